$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Test Status" column (AH) with PASS/FAIL results per row.
$ws.Range("AH1").Value = "Test Status"

$ws.Range("AH2").Value = "PASS"
$ws.Range("AH2").Interior.Color = 13434828

$ws.Range("AH3").Value = "PASS"
$ws.Range("AH3").Interior.Color = 13434828

$ws.Range("AH4").Value = "FAIL"
$ws.Range("AH4").Interior.Color = 255

$ws.Range("AH5").Value = "PASS"
$ws.Range("AH5").Interior.Color = 13434828
